$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order (A..T):
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity

$data = @(
    @("ECs","Slpi","Plscr4","ECs",1,0.3333333333333333,27.47303533333333,82.419106,0.9778496847444623,0.9778496847444623,3,1,15.281091,45.843273,0.2941125432526767,0.2941125432526767,419.8179529748819,3778.361576773937,0.287597857699022,0.2875978576990219),
    @("ECs","Slpi","Plscr4","FAPs",1,0.3333333333333333,27.47303533333333,82.419106,0.9778496847444623,0.9778496847444623,3,1,30.89280066666667,92.67840200000001,0.5945884474002099,0.5945884474002098,848.7190042609569,7638.471038348613,0.5814181258429946,0.5814181258429945),
    @("ECs","Slpi","Plscr4","sCs",1,0.3333333333333333,27.47303533333333,82.419106,0.9778496847444623,0.9778496847444623,3,1,5.782719333333333,17.348158,0.1112990093471134,0.1112990093471133,158.8688525674164,1429.819673106748,0.1088337012024458,0.1088337012024457),
    @("sCs","Slpi","Plscr4","ECs",3,1,0.622321,1.866963,0.02215031525553766,0.02215031525553766,3,1,15.281091,45.843273,0.2941125432526767,0.2941125432526767,9.509743832210999,85.58769448989899,0.006514685553654746,0.006514685553654745),
    @("sCs","Slpi","Plscr4","FAPs",3,1,0.622321,1.866963,0.02215031525553766,0.02215031525553766,3,1,30.89280066666667,92.67840200000001,0.5945884474002099,0.5945884474002098,19.22523860368067,173.027147433126,0.01317032155721532,0.01317032155721532),
    @("sCs","Slpi","Plscr4","sCs",3,1,0.622321,1.866963,0.02215031525553766,0.02215031525553766,3,1,5.782719333333333,17.348158,0.1112990093471134,0.1112990093471133,3.598707678239333,32.388369104154,0.002465308144667594,0.002465308144667594)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}
